# Fruta / hortaliza, semanal
# Insert a new weekly record as row 30, pushing the existing rows 30-86 down to 31-87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 30 (shifts rows 30..86 -> 31..87)
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with the new weekly price record
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C30").Value = 'Ñuble'
$ws.Range("D30").Value = 44614
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = 100112030
$ws.Range("G30").Value = 'Poroto granado'
$ws.Range("H30").Value = 'Sin especificar'
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 150
$ws.Range("K30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = 20000
$ws.Range("N30").Value = '$/saco 25 kilos'
$ws.Range("O30").Value = 'Provincia de Diguillín'
$ws.Range("P30").Value = 800
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = 'Hortaliza'
